# Generate Report for Handback
# Appends a new "row 4" (for the fdcad896-... file) to each of the three
# tables/worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileGuid   = "fdcad896-d91a-4005-9f55-4c8a4f009982"
$xlfHash    = "4f2da23f326d5f8ff284fb11238d556377067894"
$srcName    = "$fileGuid.md"
$srcPath    = "e2e\$fileGuid.md"
$statusText = "Handed back: in sync with en-US"

$zhXlf      = "$fileGuid.$xlfHash.zh-cn.xlf"
$deXlf      = "$fileGuid.$xlfHash.de-de.xlf"

$zhHandoffDate  = "2016-09-05 00:50:44"
$zhHandbackDate = "2016-09-05 00:51:03"
$deHandoffDate  = "2016-09-05 00:50:48"
$deHandbackDate = "2016-09-05 00:51:15"

$overviewDate = "2016-09-05 00:50:48"
$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> table3 (A1:G3 -> A1:G4)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rIdx = $loOverview.Range.Row + $loOverview.Range.Rows.Count - 1

$wsOverview.Cells.Item($rIdx, 1).Value = $srcName
$wsOverview.Cells.Item($rIdx, 2).Value = $srcPath
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($rIdx, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$xlfHash/e2e/$fileGuid.md",
    "",
    "",
    $srcPath
) | Out-Null
$wsOverview.Cells.Item($rIdx, 3).Value = ".md"
$wsOverview.Cells.Item($rIdx, 5).Value = $statusText
$wsOverview.Cells.Item($rIdx, 6).Value = $statusText
$wsOverview.Cells.Item($rIdx, 7).Value = $overviewDate
$wsOverview.Cells.Item($rIdx, 7).NumberFormat = $dateTimeFormat

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> table1 (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rIdxZh = $loZh.Range.Row + $loZh.Range.Rows.Count - 1

$wsZh.Cells.Item($rIdxZh, 1).Value = $srcName
$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rIdxZh, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$xlfHash/e2e/$fileGuid.md",
    "",
    "",
    $srcName
) | Out-Null
$wsZh.Cells.Item($rIdxZh, 2).Value = ".md"
$wsZh.Cells.Item($rIdxZh, 3).Value = $statusText
$wsZh.Cells.Item($rIdxZh, 4).Value = "e2e"
$wsZh.Cells.Item($rIdxZh, 5).Value = "ht"
$wsZh.Cells.Item($rIdxZh, 6).Value = "True"
$wsZh.Cells.Item($rIdxZh, 7).Value = $zhXlf
$wsZh.Cells.Item($rIdxZh, 8).Value = $zhHandoffDate
$wsZh.Cells.Item($rIdxZh, 8).NumberFormat = $dateTimeFormat
$wsZh.Cells.Item($rIdxZh, 9).Value = $srcName
$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rIdxZh, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$xlfHash/e2e/$fileGuid.md",
    "",
    "",
    $srcName
) | Out-Null
$wsZh.Cells.Item($rIdxZh, 10).Value = $zhXlf
$wsZh.Cells.Item($rIdxZh, 11).Value = $zhHandbackDate
$wsZh.Cells.Item($rIdxZh, 11).NumberFormat = $dateTimeFormat
$wsZh.Cells.Item($rIdxZh, 13).Value = "True"
$wsZh.Cells.Item($rIdxZh, 15).Value = "False"

# ---------------------------------------------------------------------------
# Sheet "de-de" -> table2 (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rIdxDe = $loDe.Range.Row + $loDe.Range.Rows.Count - 1

$wsDe.Cells.Item($rIdxDe, 1).Value = $srcName
$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rIdxDe, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$xlfHash/e2e/$fileGuid.md",
    "",
    "",
    $srcName
) | Out-Null
$wsDe.Cells.Item($rIdxDe, 2).Value = ".md"
$wsDe.Cells.Item($rIdxDe, 3).Value = $statusText
$wsDe.Cells.Item($rIdxDe, 4).Value = "e2e"
$wsDe.Cells.Item($rIdxDe, 5).Value = "ht"
$wsDe.Cells.Item($rIdxDe, 6).Value = "True"
$wsDe.Cells.Item($rIdxDe, 7).Value = $deXlf
$wsDe.Cells.Item($rIdxDe, 8).Value = $deHandoffDate
$wsDe.Cells.Item($rIdxDe, 8).NumberFormat = $dateTimeFormat
$wsDe.Cells.Item($rIdxDe, 9).Value = $srcName
$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rIdxDe, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$xlfHash/e2e/$fileGuid.md",
    "",
    "",
    $srcName
) | Out-Null
$wsDe.Cells.Item($rIdxDe, 10).Value = $deXlf
$wsDe.Cells.Item($rIdxDe, 11).Value = $deHandbackDate
$wsDe.Cells.Item($rIdxDe, 11).NumberFormat = $dateTimeFormat
$wsDe.Cells.Item($rIdxDe, 13).Value = "True"
$wsDe.Cells.Item($rIdxDe, 15).Value = "False"
